# Append 4 new submitted-form rows (15-18) to the "Data" sheet, mirroring
# the existing rows 8-14: columns E..M only, all stored as literal TEXT
# (mimicking the source file's numberStoredAsText / t="str" cells), and
# expand the used range (dimension / ignoredErrors) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the whole new block to be treated as Text *before* assigning any
# values, so numeric-looking strings ("345 ", "1234567890", "42567", ...)
# and the date-like "10-10-2024" strings are kept as literal text instead
# of being auto-converted to numbers / date serials by Excel.
$newRange = $ws.Range("E15:M18")
$newRange.NumberFormat = "@"

$rows = @(
    @{ Row = 15; E = "10-10-2024"; F = "345 ";   G = "P23"; H = "1234567890"; I = "854 ";  J = "client_deepak"; K = "0987654321"; L = "Submitted"; M = "Testing Reason" },
    @{ Row = 16; E = "10-10-2024"; F = "9";      G = "P23"; H = "undefined";  I = "undefined"; J = "undefined"; K = "undefined"; L = "Submitted"; M = "General Reason" },
    @{ Row = 17; E = "10-10-2024"; F = "100";    G = "P23"; H = "777888 ";    I = "33 ";    J = "seller";  K = "test009 "; L = "Submitted"; M = "Material schedule No/Service Reason" },
    @{ Row = 18; E = "10-10-2024"; F = "42567";  G = "P23"; H = "777888 ";    I = "66";     J = "gowtham"; K = "test009 "; L = "Submitted"; M = "Testing Reason" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Range("E$rowNum").Value = $r.E
    $ws.Range("F$rowNum").Value = $r.F
    $ws.Range("G$rowNum").Value = $r.G
    $ws.Range("H$rowNum").Value = $r.H
    $ws.Range("I$rowNum").Value = $r.I
    $ws.Range("J$rowNum").Value = $r.J
    $ws.Range("K$rowNum").Value = $r.K
    $ws.Range("L$rowNum").Value = $r.L
    $ws.Range("M$rowNum").Value = $r.M
}

# Best-effort: tell Excel to ignore the "number stored as text" smart-tag
# for the (now bigger) data range, so it doesn't flag the text-looking
# numbers we just inserted. xlNumberAsText = 9.
try {
    $fullRange = $ws.Range("A1:M18")
    $fullRange.Errors.Item(9).Ignore = $true
} catch {
    # Not fatal if the host doesn't support per-range error suppression.
}

Write-Host "Added rows 15-18 to sheet '$($ws.Name)'."
